$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.761.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.974.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -7.78%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.59%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.957.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -8.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.474"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -13.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -12.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -11.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.438"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -12.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -13.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000210"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -13.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.452.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.842.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.50%  "

$ws.Range("E17").Value = "  -4.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.987.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "470.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -12.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -13.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -13.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.642"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -16.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -14.07%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.65%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -11.51%  "

$ws.Range("E26").Value = "  -0.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -17.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.67%  "

$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.75%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "24.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -14.78%  "

$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.56%  "

$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "501.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -11.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -14.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0392"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0766"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -14.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.769.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.46%  "

$ws.Range("E44").Value = "  -0.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.229"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -13.16%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "115.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.38%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.104"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.43%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.97%  "

$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0497"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -14.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -12.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -18.89%  "
